$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# G1 used to read "Predict and Enter!" (a header for the "computed"/predicted
# columns); it should instead match the "Enter from gradebook:" header used
# by the other manually-entered columns (B1:F1).
$ws.Range("G1").Value = "Enter from gradebook:"

# G3 used to carry the "Note"-style orange-text formatting (reserved for the
# predicted/computed entry column); now that G is a gradebook-entry column
# like B3:F3, give it the same "Good" style those cells use.
$ws.Range("G3").Style = "Good"

# Add the new guidance text below the existing note (row 6), in rows 8-13.
$ws.Range("E8").Value = "When guessing your final project %:"
$ws.Range("E9").Value = "Middle C = 75%"
$ws.Range("E10").Value = "Middle B = 85%"
$ws.Range("E11").Value = "Middle A = 95%"
$ws.Range("E12").Value = "(To get above 95%, you'll want to go a little above and beyond the specifications, like making your python front end have a "
$ws.Range("E13").Value = "main function/menu system and/or do a really good job on comments in your queries and program.)"

# Keep the new rows on the sheet's default (custom) row height, matching the
# existing rows.
$ws.Rows.Item(8).RowHeight = 22.5
$ws.Rows.Item(9).RowHeight = 22.5
$ws.Rows.Item(10).RowHeight = 22.5
$ws.Rows.Item(11).RowHeight = 22.5
$ws.Rows.Item(12).RowHeight = 22.5
$ws.Rows.Item(13).RowHeight = 22.5

# Match the saved selection state.
[void]$ws.Range("H15").Select()
